$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Quantity" column (AI) ---
$ws.Range("AI1").Value = "Quantity"
$ws.Range("AI7").Value = 4

# --- Row 9: postcode value changed ---
$ws.Range("L9").Value = 122345

# --- Row 21: clear the old contents that are not part of the new row ---
$ws.Range("F21").ClearContents()
$ws.Range("Y21").ClearContents()
$ws.Range("Z21").ClearContents()
$ws.Range("AA21").ClearContents()
$ws.Range("AB21").ClearContents()
$ws.Range("AC21").ClearContents()

# --- Row 21: new "NewBillingAddress" test-case data ---
$ws.Range("A21").Value = "NewBillingAddress"
$ws.Range("B21").Value = "New Address"
$ws.Range("D21").Value = "harish"
$ws.Range("E21").Value = "chiruvella"
$ws.Range("G21").Value = "655 N Gabriel Ave Newton NC "
$ws.Range("H21").Value = "Lotuswave"
$ws.Range("I21").Value = "United Kingdom"
$ws.Range("J21").Value = "florida"
$ws.Range("L21").Value = 428658
$ws.Range("M21").Value = 5236987412

# --- Row 22: new "New Address" test-case data ---
$ws.Range("D22").Value = "Harish"
$ws.Range("E22").Value = "chiruvella"
$ws.Range("G22").Value = "855 Colony Dr crowley TX"
$ws.Range("H22").Value = "Lotuswave"
$ws.Range("I22").Value = "United Kingdom"
$ws.Range("J22").Value = "florida"
$ws.Range("L22").Value = 799272
$ws.Range("M22").Value = 8523697415

# --- Hyperlinks: the old F21 (email) hyperlink no longer applies; rebuild the
#     remaining 13 hyperlinks so the F21 one is dropped ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Mahendra@123.com", "", "", "Mahendra@123.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:harish.chiruvella1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:harish.chiruvella1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("U7"), "mailto:testingsstppt@gmail.com")
$ws.Hyperlinks.Add($ws.Range("V7"), "mailto:testing@123")
$ws.Hyperlinks.Add($ws.Range("F8"), "mailto:harish.chiruvella1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F9"), "mailto:harish.chiruvella1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F10"), "mailto:harish.chiruvella1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:Retailer03121CC@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F17"), "mailto:harish.chiruvella1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Retailer0112@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F19"), "mailto:harish.chiruvella1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F20"), "mailto:harish.chiruvella1@gmail.com")

# --- Selection moved to B10 ---
[void]$ws.Range("B10").Select()
